$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-CellText 'D2' '65.823.47'
Set-CellText 'E2' '  +0.49%  '
Set-CellText 'D3' '3.312.23'
Set-CellText 'E3' '  +1.57%  '
Set-CellText 'E4' '  +0.11%  '
Set-CellText 'D5' '186.78'
Set-CellText 'E5' '  +2.70%  '
Set-CellText 'D6' '552.40'
Set-CellText 'E6' '  -0.43%  '
Set-CellText 'E7' '  +0.05%  '
Set-CellText 'D8' '3.306.59'
Set-CellText 'E8' '  +1.57%  '
Set-CellText 'D9' '0.577'
Set-CellText 'E9' '  -2.29%  '
Set-CellText 'D10' '0.175'
Set-CellText 'E10' '  -5.14%  '
Set-CellText 'D11' '0.577'
Set-CellText 'E11' '  -1.42%  '
Set-CellText 'D12' '45.54'
Set-CellText 'E12' '  -3.27%  '
Set-CellText 'D13' '0.0000262'
Set-CellText 'E13' '  -0.85%  '
Set-CellText 'D14' '3.852.52'
Set-CellText 'E14' '  +1.73%  '
Set-CellText 'D15' '8.42'
Set-CellText 'E15' '  -1.51%  '
Set-CellText 'D16' '577.68'
Set-CellText 'E16' '  -8.66%  '
Set-CellText 'D17' '65.891.21'
Set-CellText 'E17' '  +0.59%  '
Set-CellText 'E18' '  +0.65%  '
Set-CellText 'D19' '3.323.76'
Set-CellText 'E19' '  +1.74%  '
Set-CellText 'D20' '17.66'
Set-CellText 'E20' '  -0.73%  '
Set-CellText 'D21' '10.82'
Set-CellText 'E21' '  -4.57%  '
Set-CellText 'D22' '0.888'
Set-CellText 'E22' '  -1.42%  '
Set-CellText 'D23' '17.88'
Set-CellText 'E23' '  +1.40%  '
Set-CellText 'D24' '4.96'
Set-CellText 'E24' '  +1.01%  '
Set-CellText 'D25' '98.20'
Set-CellText 'E25' '  -7.41%  '
Set-CellText 'D26' '3.93'
Set-CellText 'E26' '  -0.73%  '
Set-CellText 'D27' '2.66'
Set-CellText 'E27' '  +0.52%  '
Set-CellText 'D28' '9.32'
Set-CellText 'E28' '  -1.83%  '
Set-CellText 'D29' '30.38'
Set-CellText 'E29' '  +0.20%  '
Set-CellText 'D30' '8.35'
Set-CellText 'E30' '  -3.19%  '
Set-CellText 'D31' '6.57'
Set-CellText 'E31' '  +5.02%  '
Set-CellText 'D32' '571.52'
Set-CellText 'E32' '  +5.07%  '
Set-CellText 'D33' '3.68'
Set-CellText 'E33' '  -6.95%  '
Set-CellText 'D34' '10.79'
Set-CellText 'E34' '  -1.75%  '
Set-CellText 'D35' '0.102'
Set-CellText 'E35' '  -1.81%  '
Set-CellText 'D36' '3.706.35'
Set-CellText 'E36' '  +1.68%  '
Set-CellText 'E37' '  -0.25%  '
Set-CellText 'D38' '55.38'
Set-CellText 'E38' '  -2.89%  '
Set-CellText 'D39' '33.65'
Set-CellText 'E39' '  +4.42%  '
Set-CellText 'E40' '  -3.59%  '
Set-CellText 'D41' '0.0₃0683'
Set-CellText 'E41' '  -4.90%  '
Set-CellText 'B42' 'Stacks'
Set-CellText 'C42' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-CellText 'D42' '3.12'
Set-CellText 'E42' '  -7.50%  '
Set-CellText 'B43' 'ApeXProtocol'
Set-CellText 'C43' 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-CellText 'D43' '3.37'
Set-CellText 'E43' '  +2.91%  '
Set-CellText 'B44' 'Fetch.AI'
Set-CellText 'C44' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-CellText 'D44' '2.58'
Set-CellText 'E44' '  -4.66%  '
Set-CellText 'D45' '0.332'
Set-CellText 'E45' '  -0.75%  '
Set-CellText 'D46' '0.0406'
Set-CellText 'E46' '  -1.94%  '
Set-CellText 'E47' '  -0.69%  '
Set-CellText 'D48' '1.00'
Set-CellText 'E48' '  +0.16%  '
Set-CellText 'D49' '2.91'
Set-CellText 'E49' '  -13.02%  '
Set-CellText 'D50' '2.50'
Set-CellText 'E50' '  -4.19%  '
Set-CellText 'D51' '126.63'
Set-CellText 'E51' '  +5.90%  '
